$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.0258731143856077
$ws.Range("C2").Value = 0.8712523508600251
$ws.Range("D2").Value = 1.357352004991913
$ws.Range("E2").Value = 1.165054507305093
$ws.Range("F2").Value = 1.196684520570948
$ws.Range("G2").Value = 19
$ws.Range("B3").Value = -0.05889625796533703
$ws.Range("C3").Value = 0.7809309034430243
$ws.Range("D3").Value = 1.486890174830229
$ws.Range("E3").Value = 1.219381062191073
$ws.Range("F3").Value = 1.25326832255184
$ws.Range("G3").Value = 18
$ws.Range("B4").Value = -0.03036361610786248
$ws.Range("C4").Value = 0.7674964733547235
$ws.Range("D4").Value = 0.9590082094163144
$ws.Range("E4").Value = 0.9792896453125165
$ws.Range("F4").Value = 1.008943334136136
$ws.Range("G4").Value = 17
$ws.Range("B5").Value = 0.09993964655644642
$ws.Range("C5").Value = 0.6259185880536732
$ws.Range("D5").Value = 0.8952003656645542
$ws.Range("E5").Value = 0.9461502870393024
$ws.Range("F5").Value = 0.9717132266730971
$ws.Range("G5").Value = 16
$ws.Range("B6").Value = 0.03831129326750842
$ws.Range("C6").Value = 0.7119710808408777
$ws.Range("D6").Value = 1.080030924484957
$ws.Range("E6").Value = 1.039245362984583
$ws.Range("F6").Value = 1.074989951428016
$ws.Range("G6").Value = 15
$ws.Range("B7").Value = 0.09938456961952043
$ws.Range("C7").Value = 0.7378794059141224
$ws.Range("D7").Value = 1.55852155670565
$ws.Range("E7").Value = 1.248407608397854
$ws.Range("F7").Value = 1.291421985981127
$ws.Range("G7").Value = 14
$ws.Range("B8").Value = 0.06981294092732569
$ws.Range("C8").Value = 0.5769263726184313
$ws.Range("D8").Value = 0.6329722743072022
$ws.Range("E8").Value = 0.7955955469377655
$ws.Range("F8").Value = 0.8248878488730877
$ws.Range("G8").Value = 13
$ws.Range("B9").Value = 0.2198501169052137
$ws.Range("C9").Value = 0.8740014322694108
$ws.Range("D9").Value = 1.535812380474684
$ws.Range("E9").Value = 1.239278976048042
$ws.Range("F9").Value = 1.273853840583311
$ws.Range("G9").Value = 12
$ws.Range("B10").Value = 0.1257302919151607
$ws.Range("C10").Value = 0.7877651198835861
$ws.Range("D10").Value = 0.9791626505880603
$ws.Range("E10").Value = 0.9895264779620908
$ws.Range("F10").Value = 1.029412453155336
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = 0.2215020228704644
$ws.Range("C11").Value = 0.8099270631829588
$ws.Range("D11").Value = 1.114429208763017
$ws.Range("E11").Value = 1.055665292014006
$ws.Range("F11").Value = 1.08799819374202
$ws.Range("G11").Value = 10
